$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: reduce numeric precision (2 decimal places) to match the
# "custom accuracy" export. B5, Z5 and AA5 are unchanged.
$ws.Range("C5").Value  = 7.45
$ws.Range("D5").Value  = 0.76
$ws.Range("E5").Value  = 21.88
$ws.Range("F5").Value  = 17.68
$ws.Range("G5").Value  = 7.5
$ws.Range("H5").Value  = 33.91
$ws.Range("I5").Value  = 12.09
$ws.Range("J5").Value  = 5.73
$ws.Range("K5").Value  = 7.67
$ws.Range("L5").Value  = 8.800000000000001
$ws.Range("M5").Value  = 9.43
$ws.Range("N5").Value  = 2.72
$ws.Range("O5").Value  = 7.9
$ws.Range("P5").Value  = 11.03
$ws.Range("Q5").Value  = 6.81
$ws.Range("R5").Value  = 0.22
$ws.Range("S5").Value  = 0.47
$ws.Range("T5").Value  = 112.86
$ws.Range("U5").Value  = 22.2
$ws.Range("V5").Value  = 7.29
$ws.Range("W5").Value  = 14.74
$ws.Range("X5").Value  = 7.9
$ws.Range("Y5").Value  = 1.04
$ws.Range("AB5").Value = 5.78
$ws.Range("AC5").Value = 6.85
$ws.Range("AD5").Value = 9.300000000000001
$ws.Range("AE5").Value = 0.47
$ws.Range("AF5").Value = 31.09
$ws.Range("AG5").Value = 4.12
$ws.Range("AH5").Value = 9.050000000000001

# Drop the now-superfluous last data row (row 6) entirely; this also
# shrinks the sheet's used range / dimension down to A1:AH5.
$ws.Rows.Item(6).Delete()
